$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1846.8823
$ws.Range("I86").Value = 2063.125
$ws.Range("J86").Value = 1654.6666
$ws.Range("K86").Value = 2063.125
$ws.Range("L86").Value = 1654.6666
$ws.Range("M86").Value = -940.125
$ws.Range("N86").Value = -3900.6666

$ws.Range("H89").Value = 1846.8823
$ws.Range("I89").Value = 2063.125
$ws.Range("J89").Value = 1654.6666
$ws.Range("K89").Value = 10315.625
$ws.Range("L89").Value = 8273.333000000001
$ws.Range("M89").Value = -4699.625
$ws.Range("N89").Value = -19505.333

$ws.Range("H92").Value = 768.1429000000001
$ws.Range("I92").Value = 752.73334
$ws.Range("K92").Value = 752.73334
$ws.Range("M92").Value = 495.26666

$ws.Range("H97").Value = 72441.14
$ws.Range("I97").Value = 100
$ws.Range("J97").Value = 84498
$ws.Range("K97").Value = 300
$ws.Range("L97").Value = 253494
$ws.Range("M97").Value = 196
$ws.Range("N97").Value = -254486

$ws.Range("H112").Value = 1902.5333
$ws.Range("J112").Value = 1902.5333
$ws.Range("L112").Value = 5707.5999
$ws.Range("N112").Value = -7923.5999

$ws.Range("H129").Value = 2915.4695
$ws.Range("I129").Value = 25563.25
$ws.Range("J129").Value = 902.3333
$ws.Range("K129").Value = 76689.75
$ws.Range("L129").Value = 2706.9999
$ws.Range("M129").Value = -71689.75
$ws.Range("N129").Value = -12706.9999

$ws.Range("H132").Value = 4172059.8
$ws.Range("I132").Value = 4634992.5
$ws.Range("J132").Value = 5666.5
$ws.Range("K132").Value = 13904977.5
$ws.Range("L132").Value = 16999.5
$ws.Range("M132").Value = -13902447.5
$ws.Range("N132").Value = -22059.5

$ws.Range("H138").Value = 2453.3
$ws.Range("I138").Value = 1136.9546
$ws.Range("J138").Value = 2824.577
$ws.Range("K138").Value = 3410.8638
$ws.Range("L138").Value = 8473.731
$ws.Range("M138").Value = 1729.1362
$ws.Range("N138").Value = -18753.731

$ws.Range("H141").Value = 2724.7
$ws.Range("I141").Value = 2138
$ws.Range("J141").Value = 4484.8
$ws.Range("K141").Value = 6414
$ws.Range("L141").Value = 13454.4
$ws.Range("M141").Value = -1234
$ws.Range("N141").Value = -23814.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26919.771
$ws.Range("I32").Value = 8754.556
$ws.Range("J32").Value = 145819.36
$ws.Range("K32").Value = 8754.556
$ws.Range("L32").Value = 145819.36
$ws.Range("M32").Value = -8467.556
$ws.Range("N32").Value = -146393.36

$ws.Range("H61").Value = 1688.963
$ws.Range("I61").Value = 1220.7142
$ws.Range("J61").Value = 2193.2307
$ws.Range("K61").Value = 1220.7142
$ws.Range("L61").Value = 2193.2307
$ws.Range("M61").Value = -1008.7142
$ws.Range("N61").Value = -2617.2307

$ws.Range("H74").Value = 972.675
$ws.Range("I74").Value = 916.4483
$ws.Range("K74").Value = 916.4483
$ws.Range("M74").Value = -42.44830000000002

$ws.Range("H77").Value = 972.675
$ws.Range("I77").Value = 916.4483
$ws.Range("K77").Value = 4582.2415
$ws.Range("M77").Value = -214.2415000000001

$ws.Range("H122").Value = 2161.9678
$ws.Range("I122").Value = 1782.3182
$ws.Range("K122").Value = 5346.9546
$ws.Range("M122").Value = -2896.9546

$ws.Range("H132").Value = 25158.371
$ws.Range("I132").Value = 40252.438
$ws.Range("J132").Value = 3203.3635
$ws.Range("K132").Value = 120757.314
$ws.Range("L132").Value = 9610.0905
$ws.Range("M132").Value = -118227.314
$ws.Range("N132").Value = -14670.0905

$ws.Range("H136").Value = 1688.963
$ws.Range("I136").Value = 1220.7142
$ws.Range("J136").Value = 2193.2307
$ws.Range("K136").Value = 3662.1426
$ws.Range("L136").Value = 6579.6921
$ws.Range("M136").Value = -1112.1426
$ws.Range("N136").Value = -11679.6921

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("N108").Value = 0

$ws.Range("H134").Value = 2919.7273
$ws.Range("I134").Value = 2706.5386
$ws.Range("J134").Value = 4582.6
$ws.Range("K134").Value = 8119.6158
$ws.Range("L134").Value = 13747.8
$ws.Range("M134").Value = -5584.6158
$ws.Range("N134").Value = -18817.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 38086.574
$ws.Range("I31").Value = 1113.0714
$ws.Range("J31").Value = 57995.383
$ws.Range("K31").Value = 1113.0714
$ws.Range("L31").Value = 57995.383
$ws.Range("M31").Value = -818.0714
$ws.Range("N31").Value = -58585.383

$ws.Range("H34").Value = 38086.574
$ws.Range("I34").Value = 1113.0714
$ws.Range("J34").Value = 57995.383
$ws.Range("K34").Value = 1113.0714
$ws.Range("L34").Value = 57995.383
$ws.Range("M34").Value = -911.0714
$ws.Range("N34").Value = -58399.383

$ws.Range("H58").Value = 1609.1818
$ws.Range("I58").Value = 1465.7307
$ws.Range("J58").Value = 2142
$ws.Range("K58").Value = 1465.7307
$ws.Range("L58").Value = 2142
$ws.Range("M58").Value = -1262.7307
$ws.Range("N58").Value = -2548

$ws.Range("H62").Value = 2300.3333
$ws.Range("I62").Value = 1500
$ws.Range("J62").Value = 2460.4
$ws.Range("K62").Value = 1500
$ws.Range("L62").Value = 2460.4
$ws.Range("M62").Value = -876
$ws.Range("N62").Value = -3708.4

$ws.Range("H65").Value = 2300.3333
$ws.Range("I65").Value = 1500
$ws.Range("J65").Value = 2460.4
$ws.Range("K65").Value = 7500
$ws.Range("L65").Value = 12302
$ws.Range("M65").Value = -4380
$ws.Range("N65").Value = -18542

$ws.Range("H132").Value = 2350.7568
$ws.Range("I132").Value = 2205.2942
$ws.Range("J132").Value = 3999.3333
$ws.Range("K132").Value = 6615.882599999999
$ws.Range("L132").Value = 11997.9999
$ws.Range("M132").Value = -4085.882599999999
$ws.Range("N132").Value = -17057.9999

$ws.Range("H136").Value = 1609.1818
$ws.Range("I136").Value = 1465.7307
$ws.Range("J136").Value = 2142
$ws.Range("K136").Value = 4397.1921
$ws.Range("L136").Value = 6426
$ws.Range("M136").Value = -1847.1921
$ws.Range("N136").Value = -11526

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 113744.664
$ws.Range("I70").Value = 335267.34
$ws.Range("J70").Value = 2983.3333
$ws.Range("K70").Value = 1005802.02
$ws.Range("L70").Value = 8949.999899999999
$ws.Range("M70").Value = -1005487.02
$ws.Range("N70").Value = -9579.999899999999

$ws.Range("H73").Value = 113744.664
$ws.Range("I73").Value = 335267.34
$ws.Range("J73").Value = 2983.3333
$ws.Range("K73").Value = 1005802.02
$ws.Range("L73").Value = 8949.999899999999
$ws.Range("M73").Value = -1004710.02
$ws.Range("N73").Value = -11133.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 91111280
$ws.Range("I80").Value = 250551250
$ws.Range("J80").Value = 2725.7144
$ws.Range("K80").Value = 250551250
$ws.Range("L80").Value = 2725.7144
$ws.Range("M80").Value = -250550252
$ws.Range("N80").Value = -4721.7144

$ws.Range("H83").Value = 91111280
$ws.Range("I83").Value = 250551250
$ws.Range("J83").Value = 2725.7144
$ws.Range("K83").Value = 1252756250
$ws.Range("L83").Value = 13628.572
$ws.Range("M83").Value = -1252751258
$ws.Range("N83").Value = -23612.572

$ws.Range("H126").Value = 3270976
$ws.Range("I126").Value = 3361.0908
$ws.Range("J126").Value = 8405799
$ws.Range("K126").Value = 10083.2724
$ws.Range("L126").Value = 25217397
$ws.Range("M126").Value = -7613.2724
$ws.Range("N126").Value = -25222337

$ws.Range("H132").Value = 6769.7
$ws.Range("I132").Value = 6175
$ws.Range("K132").Value = 18525
$ws.Range("M132").Value = -15995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3852.4546
$ws.Range("I68").Value = 2223.3845
$ws.Range("J68").Value = 6205.5557
$ws.Range("K68").Value = 2223.3845
$ws.Range("L68").Value = 6205.5557
$ws.Range("M68").Value = -1474.3845
$ws.Range("N68").Value = -7703.5557

$ws.Range("H71").Value = 3852.4546
$ws.Range("I71").Value = 2223.3845
$ws.Range("J71").Value = 6205.5557
$ws.Range("K71").Value = 11116.9225
$ws.Range("L71").Value = 31027.7785
$ws.Range("M71").Value = -7372.922500000001
$ws.Range("N71").Value = -38515.7785

$ws.Range("H122").Value = 3149.5386
$ws.Range("I122").Value = 3028.2222
$ws.Range("J122").Value = 3422.5
$ws.Range("K122").Value = 9084.6666
$ws.Range("L122").Value = 10267.5
$ws.Range("M122").Value = -6634.6666
$ws.Range("N122").Value = -15167.5

$ws.Range("H132").Value = 4399.095
$ws.Range("I132").Value = 4970.5713
$ws.Range("K132").Value = 14911.7139
$ws.Range("M132").Value = -12381.7139

$ws.Range("H136").Value = 2660
$ws.Range("I136").Value = 2008.2858
$ws.Range("J136").Value = 3420.3333
$ws.Range("K136").Value = 6024.857400000001
$ws.Range("L136").Value = 10260.9999
$ws.Range("M136").Value = -3474.857400000001
$ws.Range("N136").Value = -15360.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 83879.414
$ws.Range("I107").Value = 625.6667
$ws.Range("J107").Value = 167133.17
$ws.Range("K107").Value = 1877.0001
$ws.Range("L107").Value = 501399.51
$ws.Range("M107").Value = 42.99990000000003
$ws.Range("N107").Value = -505239.51

$ws.Range("H132").Value = 6706.467
$ws.Range("I132").Value = 5636.4546
$ws.Range("J132").Value = 7325.9473
$ws.Range("K132").Value = 16909.3638
$ws.Range("L132").Value = 21977.8419
$ws.Range("M132").Value = -14379.3638
$ws.Range("N132").Value = -27037.8419

$ws.Range("H136").Value = 21116.965
$ws.Range("I136").Value = 50819.1
$ws.Range("J136").Value = 5061.757
$ws.Range("K136").Value = 152457.3
$ws.Range("L136").Value = 15185.271
$ws.Range("M136").Value = -149907.3
$ws.Range("N136").Value = -20285.271
